$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.406.68'
$ws.Range('E2').Value = '  +2.32%  '

# Row 3
$ws.Range('D3').Value = '2.095.33'
$ws.Range('E3').Value = '  -0.04%  '

# Row 4
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = '  -0.74%  '

# Row 5
$ws.Range('D5').Value = "'343.17"
$ws.Range('E5').Value = '  -0.08%  '

# Row 6
$ws.Range('E6').Value = '  -0.66%  '

# Row 7
$ws.Range('D7').Value = "'0.5252"
$ws.Range('E7').Value = '  +1.67%  '

# Row 8
$ws.Range('D8').Value = "'0.4428"
$ws.Range('E8').Value = '  +1.06%  '

# Row 9
$ws.Range('D9').Value = "'54.54"
$ws.Range('E9').Value = '  +3.38%  '

# Row 10
$ws.Range('D10').Value = "'0.09366"
$ws.Range('E10').Value = '  +0.92%  '

# Row 11
$ws.Range('D11').Value = "'1.171"
$ws.Range('E11').Value = '  +0.54%  '

# Row 12
$ws.Range('E12').Value = '  -0.16%  '

# Row 13
$ws.Range('D13').Value = '2.149.56'
$ws.Range('E13').Value = '  +2.49%  '

# Row 14
$ws.Range('D14').Value = "'8.634"
$ws.Range('E14').Value = '  +4.11%  '

# Row 15
$ws.Range('D15').Value = "'6.916"
$ws.Range('E15').Value = '  +2.46%  '

# Row 16
$ws.Range('D16').Value = "'101.70"
$ws.Range('E16').Value = '  +2.23%  '

# Row 17
$ws.Range('D17').Value = "'0.00001161"
$ws.Range('E17').Value = '  +0.78%  '

# Row 18
$ws.Range('E18').Value = '  -0.62%  '

# Row 19
$ws.Range('D19').Value = "'21.20"
$ws.Range('E19').Value = '  +2.00%  '

# Row 20
$ws.Range('D20').Value = "'0.06684"
$ws.Range('E20').Value = '  +0.39%  '

# Row 21
$ws.Range('D21').Value = "'6.344"
$ws.Range('E21').Value = '  +2.36%  '

# Row 22
$ws.Range('D22').Value = "'1.001"
$ws.Range('E22').Value = '  -0.68%  '

# Row 23
$ws.Range('D23').Value = '30.412.09'

# Row 24
$ws.Range('D24').Value = "'12.56"
$ws.Range('E24').Value = '  +0.27%  '

# Row 25
$ws.Range('D25').Value = "'2.312"
$ws.Range('E25').Value = '  -0.21%  '

# Row 26
$ws.Range('D26').Value = "'21.88"
$ws.Range('E26').Value = '  -0.33%  '

# Row 27
$ws.Range('D27').Value = "'162.93"
$ws.Range('E27').Value = '  +1.09%  '

# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = "'6.828"
$ws.Range('E28').Value = '  +9.02%  '

# Row 29
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = "'2.516"
$ws.Range('E29').Value = '  -0.21%  '

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = "'133.68"
$ws.Range('E30').Value = '  +0.36%  '

# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = "'1.141"
$ws.Range('E31').Value = '  +0.39%  '

# Row 32
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = "'1.664"
$ws.Range('E32').Value = '  +0.68%  '

# Row 33
$ws.Range('B33').Value = 'Stellar'
$ws.Range('C33').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D33').Value = "'0.1049"
$ws.Range('E33').Value = '  -0.08%  '

# Row 34
$ws.Range('D34').Value = "'6.278"
$ws.Range('E34').Value = '  +1.78%  '

# Row 35
$ws.Range('D35').Value = "'3.866"
$ws.Range('E35').Value = '  -1.82%  '

# Row 36
$ws.Range('D36').Value = "'10.18"
$ws.Range('E36').Value = '  -0.47%  '

# Row 37
$ws.Range('D37').Value = "'0.02638"

# Row 38
$ws.Range('D38').Value = "'0.06832"
$ws.Range('E38').Value = '  +1.78%  '

# Row 39
$ws.Range('D39').Value = "'0.7002"
$ws.Range('E39').Value = '  +1.61%  '

# Row 40
$ws.Range('D40').Value = "'12.62"
$ws.Range('E40').Value = '  +1.01%  '

# Row 41
$ws.Range('E41').Value = '  +1.82%  '

# Row 42
$ws.Range('D42').Value = "'0.2223"
$ws.Range('E42').Value = '  +0.07%  '

# Row 43
$ws.Range('D43').Value = "'0.6830"
$ws.Range('E43').Value = '  +0.82%  '

# Row 44
$ws.Range('D44').Value = "'14.44"
$ws.Range('E44').Value = '  +1.10%  '

# Row 45
$ws.Range('D45').Value = "'2.349"
$ws.Range('E45').Value = '  +1.20%  '

# Row 46
$ws.Range('E46').Value = '  -0.65%  '

# Row 47
$ws.Range('D47').Value = "'1.386"
$ws.Range('E47').Value = '  +19.38%  '

# Row 48
$ws.Range('D48').Value = "'3.638"
$ws.Range('E48').Value = '  +0.67%  '

# Row 49
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = "'0.00000000353"
$ws.Range('E49').Value = '  -0.41%  '

# Row 50
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = "'1.240"
$ws.Range('E50').Value = '  +10.60%  '

# Row 51
$ws.Range('D51').Value = "'1.218"
$ws.Range('E51').Value = '  -0.08%  '
